# Updated symbol list (coin prices/ranks + a block reorder of rows 9-17,
# where "One" moves up from row 17 to row 9 and the coins that were in
# rows 9-16 each shift down one row).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) cells store numeric-looking text (e.g. "245.10"); force
# text formatting first so Excel keeps them as strings instead of coercing
# to numbers (which would also lose trailing zeros).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '245.10'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '23.04'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.415'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.06003'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '3.390'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8084'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.9264'
$ws.Range("B9").Value = 'One'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.01118'
$ws.Range("E9").Value = '8OneONEBestin24h'
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1417'
$ws.Range("E10").Value = '9WazirXWRX'
$ws.Range("B11").Value = 'MandalaExchangeToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07443'
$ws.Range("E11").Value = '10MandalaExchangeTokenMDX'
$ws.Range("B12").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C12").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.03376'
$ws.Range("E12").Value = '11LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.03029'
$ws.Range("E13").Value = '12BitrueCoinBTR'
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.09354'
$ws.Range("E14").Value = '13BitMartTokenBMX'
$ws.Range("B15").Value = 'MCDex'
$ws.Range("C15").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.942'
$ws.Range("E15").Value = '14MCDexMCB'
$ws.Range("B16").Value = 'BitForexToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.001588'
$ws.Range("E16").Value = '15BitForexTokenBF'
$ws.Range("B17").Value = 'CoinExToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.04789'
$ws.Range("E17").Value = '16CoinExTokenCET'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.005418'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.004158'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0009848'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.00007702'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '6.443'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.03953'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006202'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.006662'
$ws.Range("E44").Value = '43LocalTradersLCT'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005199'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.100'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.002026'
$ws.Range("E49").Value = '48BOLOBOLOWorstin24h'
